$d = $word.ActiveDocument

# There are four paragraphs in the body that contain the old (French/English
# mixed) "Waktu Kampanye 2018 untuk Perseus: ..." text, built out of two or
# three separate runs (sometimes preceded by a "www.globeatnight.org" run
# ending in a line break). Each one must become a single, completely
# unformatted run reading "Waktu Kampanye Leo: 14-23 April, 14-23 Mei".
#
# Because replacing a paragraph's content changes character offsets for
# every paragraph that follows it, we repeatedly re-scan the document for
# the next remaining match and fix it, instead of computing all offsets up
# front.

$found = $true
while ($found) {
    $found = $false
    foreach ($p in $d.Paragraphs) {
        if ($p.Range.Text -like "*untuk Perseus*") {
            $start = $p.Range.Start
            $end = $p.Range.End

            # Wipe out every run in the paragraph (but keep the paragraph
            # mark itself, i.e. stop one character short of $end).
            $body = $d.Range($start, $end - 1)
            $body.Text = ""

            # Insert the replacement text into the now-empty paragraph. The
            # resulting run picks up no formatting at all, matching the
            # target markup (a bare <w:r><w:t>...</w:t></w:r>).
            $insertionPoint = $d.Range($start, $start)
            $insertionPoint.InsertAfter("Waktu Kampanye Leo: 14-23 April, 14-23 Mei")

            $found = $true
            break
        }
    }
}
